$wb = $excel.ActiveWorkbook

# --- Update Hoja1!A1 conversion note text ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$nl = [char]10
$newText = "Conversión del día 💰" + $nl + `
"✅ Dólar paralelo: 68" + $nl + `
"" + $nl + `
"Binance" + $nl + `
"✅ 1000 Bs = 5.45 = 21701.32 pesos" + $nl + `
"✅ 21701.32 pesos = 5.43 = 960.89 Bs" + $nl + `
"" + $nl + `
"Promedio competencia" + $nl + `
"✅ Tasa pesos: 20" + $nl + `
"✅ Tasa Bs: 20" + $nl + `
"✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $newText

# --- Update tasas sheet rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 183.399
$wsTasas.Range("O10").Value = 3980
$wsTasas.Range("N12").Value = 3999.95
$wsTasas.Range("O12").Value = 177.11
